# Append: 2026-01-22 18:29 JST
# Update the "取得日時" (retrieved-at) timestamp in column A for all data rows
# on the active sheet ("ランサーズ"), and refresh the details for the
# MySQL/MariaDB job row (row 6): its price range / terms text and priority score.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2026-01-22 18:29:40"

# Rows 2-8 hold the scraped listings; update their "取得日時" column (A).
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Row 6 ("【急募】MySQL/MariaDBを活用したデータベース設計・運用の依頼") also had its
# price/terms (D) and priority score (G) refreshed on this re-scrape.
$ws.Range("D6").Value = "10,000 円 ~ 20,000 円 / 募集期間 5 日、取引期間 0 日"
$ws.Range("G6").Value = 30
